# Fixed the VSS fail testcase
# Row 19 (VSS-18 "Test if backup/restore gracefully fails in case of failure.")
# gets an updated Scripts (F19) and XML (G19) cell: the test now adds a hard
# disk via setup/cleanup scripts and passes SCSI/IDE/FILESYS params.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VSS")

$ws.Range("F19").Value = "AddHardDisk.ps1`nVSS_BackupRestore_Fail.ps1`nVSS_Disk_Fail.sh"

$ws.Range("G19").Value = "    <test>`n    <testName>VSS_BackupRestore_Fail</testName>`n        <setupScript>setupscripts\AddHardDisk.ps1</setupScript>`n        <testScript>setupscripts\VSS_BackupRestore_Fail.ps1</testScript> `n        <testParams>`n            <param>driveletter=F:</param>`n            <param>SCSI=0,1,Dynamic</param>`n            <param>IDE=0,1,Dynamic</param>`n            <param>FILESYS-ext3</param>`n            <param>TC_COVERED=VSS-18</param>`n        </testParams>`n        <cleanupScript>setupscripts\RemoveHardDisk.ps1</cleanupScript>`n        <timeout>1200</timeout>`n        <OnERROR>Continue</OnERROR>`n    </test>"

$ws.Range("C19").Select()
